# "Generate Report for Handoff" - refresh the localization-status report:
#  - bump the two "Latest ... Xliff Generate Date" / "Latest Handoff Datetime"
#    timestamps that correspond to this handoff run
#  - set the Priority column to "ht" (handoff type) for every row whose
#    Status is "Ready for handoff", on both locale sheets

$wb = $excel.ActiveWorkbook

$rows = @(7, 8, 10, 11, 13, 14)

# Overview sheet: column G = "Latest HO Xliff Generate Date"
$overview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $overview.Cells.Item($r, 7).Value = "2016-08-16 08:21:33"
}

# zh-cn sheet: column H = "Latest Handoff Datetime", column E = "Priority"
$zhcn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $zhcn.Cells.Item($r, 8).Value = "2016-08-16 08:21:28"
    $zhcn.Cells.Item($r, 5).Value = "ht"
}

# de-de sheet: column E = "Priority"
$dede = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $dede.Cells.Item($r, 5).Value = "ht"
}
